$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.782.03"
$ws.Range("E2").Value = "  -4.74%  "

# Row 3
$ws.Range("D3").Value = "2.209.92"
$ws.Range("E3").Value = "  -5.97%  "

# Row 4
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").Value = "246.42"
$ws.Range("E5").Value = "  +2.62%  "

# Row 6
$ws.Range("E6").Value = "  -6.31%  "

# Row 7
$ws.Range("D7").Value = "69.94"
$ws.Range("E7").Value = "  -5.29%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  -7.11%  "

# Row 10
$ws.Range("D10").Value = "0.0954"
$ws.Range("E10").Value = "  -5.62%  "

# Row 11
$ws.Range("D11").Value = "57.96"
$ws.Range("E11").Value = "  -6.06%  "

# Row 12
$ws.Range("D12").Value = "35.91"
$ws.Range("E12").Value = "  +7.88%  "

# Row 13
$ws.Range("E13").Value = "  -3.24%  "

# Row 14
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -7.21%  "

# Row 15
$ws.Range("D15").Value = "2.529.50"
$ws.Range("E15").Value = "  -6.47%  "

# Row 16
$ws.Range("D16").Value = "14.86"
$ws.Range("E16").Value = "  -8.24%  "

# Row 17
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  -6.77%  "

# Row 18
$ws.Range("D18").Value = "2.194.96"
$ws.Range("E18").Value = "  -6.52%  "

# Row 19
$ws.Range("D19").Value = "41.633.89"
$ws.Range("E19").Value = "  -4.92%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  -6.83%  "

# Row 21
$ws.Range("D21").Value = "72.90"
$ws.Range("E21").Value = "  -6.45%  "

# Row 22
$ws.Range("D22").Value = "6.08"
$ws.Range("E22").Value = "  -8.07%  "

# Row 23
$ws.Range("D23").Value = "234.90"
$ws.Range("E23").Value = "  -6.96%  "

# Row 24
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  +10.83%  "

# Row 25
$ws.Range("E25").Value = "  -0.09%  "

# Row 26
$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.61"
$ws.Range("E26").Value = "  -5.53%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -1.53%  "

# Row 28
$ws.Range("D28").Value = "2.25"
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("D29").Value = "9.90"
$ws.Range("E29").Value = "  -4.91%  "

# Row 30
$ws.Range("D30").Value = "169.70"
$ws.Range("E30").Value = "  -3.36%  "

# Row 31
$ws.Range("D31").Value = "20.45"
$ws.Range("E31").Value = "  -8.01%  "

# Row 32
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -6.15%  "

# Row 33
$ws.Range("E33").Value = "  -7.41%  "

# Row 34
$ws.Range("D34").Value = "0.0716"
$ws.Range("E34").Value = "  -3.90%  "

# Row 35
$ws.Range("D35").Value = "5.14"
$ws.Range("E35").Value = "  -4.18%  "

# Row 36
$ws.Range("D36").Value = "4.66"
$ws.Range("E36").Value = "  -7.69%  "

# Row 37
$ws.Range("D37").Value = "3.85"
$ws.Range("E37").Value = "  +1.82%  "

# Row 38
$ws.Range("D38").Value = "22.85"
$ws.Range("E38").Value = "  +18.32%  "

# Row 39
$ws.Range("D39").Value = "2.29"
$ws.Range("E39").Value = "  -5.01%  "

# Row 40
$ws.Range("D40").Value = "0.0271"
$ws.Range("E40").Value = "  -0.39%  "

# Row 41
$ws.Range("D41").Value = "5.85"
$ws.Range("E41").Value = "  -9.11%  "

# Row 42
$ws.Range("D42").Value = "65.71"
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  -11.46%  "

# Row 44
$ws.Range("D44").Value = "8.91"
$ws.Range("E44").Value = "  -2.59%  "

# Row 45
$ws.Range("E45").Value = "  -4.95%  "

# Row 46
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.190"
$ws.Range("E46").Value = "  -5.10%  "

# Row 47
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.19%  "

# Row 48
$ws.Range("D48").Value = "4.59"
$ws.Range("E48").Value = "  +9.10%  "

# Row 49
$ws.Range("D49").Value = "10.23"
$ws.Range("E49").Value = "  +7.81%  "

# Row 50
$ws.Range("D50").Value = "1.18"
$ws.Range("E50").Value = "  -3.47%  "

# Row 51
$ws.Range("E51").Value = "  -4.20%  "
